$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.456.03'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '1.573.89'
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.64'
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3730'
$ws.Range("E7").Value = '  -0.78%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.95'
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3403'
$ws.Range("E9").Value = '  -0.53%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07572'
$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.145'
$ws.Range("E11").Value = '  -0.43%  '

$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.33'
$ws.Range("E13").Value = '  +0.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.028'
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.973'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("D16").Value = '1.572.69'
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001123'
$ws.Range("E17").Value = '  -0.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.98'
$ws.Range("E18").Value = '  +1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06753'
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.303'
$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("E22").Value = '  -2.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.19'
$ws.Range("E23").Value = '  +1.65%  '

$ws.Range("D24").Value = '22.450.58'
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.341'
$ws.Range("E25").Value = '  -2.31%  '

$ws.Range("E26").Value = '  +1.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.09'
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.62'
$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.013'
$ws.Range("E29").Value = '  -0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.65'

$ws.Range("D31").Value = '1.748.15'
$ws.Range("E31").Value = '  +0.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.057'
$ws.Range("E32").Value = '  +7.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.205'
$ws.Range("E33").Value = '  +0.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.985'
$ws.Range("E34").Value = '  -1.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.855'
$ws.Range("E35").Value = '  -0.39%  '

$ws.Range("E36").Value = '  -1.16%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02498'
$ws.Range("E37").Value = '  -1.80%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.356'
$ws.Range("E38").Value = '  -0.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2302'
$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06518'
$ws.Range("E40").Value = '  -0.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.482'
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.31'
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6225'
$ws.Range("E43").Value = '  -2.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.99'
$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.812'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5815'
$ws.Range("E47").Value = '  -2.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.63'
$ws.Range("E48").Value = '  +3.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.074'
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("E50").Value = '  -5.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07326'
$ws.Range("E51").Value = '  +0.03%  '
